$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A item names (hyphen instead of space) ---
$ws.Range("A2").Value = "BARANG-1"
$ws.Range("A3").Value = "BARANG-2"
$ws.Range("A4").Value = "BARANG-3"

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "NAMABARANG"
$ws.Range("W1").Value = "NAMASUPPLIER"

# --- Row 3: clear out the second unit (SATUAN2) block ---
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()

# --- Sheet view: scroll position & selection ---
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("X9").Select()

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 13.6
$ws.Columns.Item(23).ColumnWidth = 14.1
